$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update price column (B) for rows 49-52 from "450 Tl" to "350 Tl"
$ws.Range("B49:B52").Value = "350 Tl"

# Reflect the last active selection on the sheet (B49)
$ws.Range("B49").Select()
